$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.359
$ws.Range("D7").Value = -7.15
$ws.Range("C9").Value = -10.612
$ws.Range("D12").Value = -7.697999999999999
$ws.Range("C13").Value = -12.867
$ws.Range("D14").Value = -7.959999999999999
$ws.Range("E15").Value = 16.16699999999999
$ws.Range("C16").Value = -12.796
$ws.Range("C18").Value = -11.759
$ws.Range("D19").Value = -8.1
$ws.Range("C20").Value = -12.578
$ws.Range("C26").Value = -12.983
$ws.Range("D26").Value = -7.875
$ws.Range("C27").Value = -12.855
$ws.Range("D27").Value = -7.988999999999999
$ws.Range("E28").Value = 16.854
$ws.Range("C29").Value = -12.74
$ws.Range("D29").Value = -7.49
$ws.Range("E33").Value = 16.832
$ws.Range("C35").Value = -12.383
$ws.Range("E35").Value = 16.451
$ws.Range("C36").Value = -12.428
$ws.Range("D37").Value = -8.059999999999999
$ws.Range("D38").Value = -7.628
$ws.Range("E38").Value = 16.697
$ws.Range("E43").Value = 16.962
$ws.Range("E44").Value = 16.799
$ws.Range("C45").Value = -13.3
$ws.Range("E45").Value = 16.667
$ws.Range("D47").Value = -7.467999999999999
$ws.Range("E47").Value = 16.593
$ws.Range("D51").Value = -8.149000000000001
$ws.Range("E51").Value = 16.724
$ws.Range("D52").Value = -7.600999999999999
$ws.Range("E54").Value = 16.335
$ws.Range("C55").Value = -13.44
$ws.Range("C57").Value = -13.174
$ws.Range("E57").Value = 16.569
$ws.Range("E62").Value = 16.282
$ws.Range("E63").Value = 17.458
$ws.Range("E67").Value = 17.074
$ws.Range("C69").Value = -11.04
$ws.Range("D69").Value = -7.178
$ws.Range("D70").Value = -7.434
$ws.Range("E70").Value = 17.421
$ws.Range("C76").Value = -13.42
$ws.Range("D76").Value = -7.453999999999999
$ws.Range("C78").Value = -12.577
$ws.Range("D81").Value = -8.403
$ws.Range("E81").Value = 16.602
$ws.Range("C82").Value = -11.884
$ws.Range("C83").Value = -13.184
$ws.Range("D83").Value = -8.472999999999999
$ws.Range("E88").Value = 16.488
$ws.Range("C93").Value = -11.704
$ws.Range("D94").Value = -7.455999999999999
$ws.Range("E96").Value = 16.482
$ws.Range("C97").Value = -12.35
$ws.Range("E99").Value = 16.632
$ws.Range("D100").Value = -8.221
$ws.Range("D102").Value = -7.764
